# Actualización automática 2025-09-11 14:30:08
$wb = $excel.ActiveWorkbook

# Sheet 1: VENTAS POR GRUPO
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M14").Value = 11342.53

# Sheet 2: VENTA MENSUAL
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F14").Value = 11342.53
$ws2.Range("F23").Value = 31430.95

# Sheet 3: CUMPLIMIENTO MENSUAL
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D12").Value = 29337.7
$ws3.Range("E12").Value = 7485.943092117097
$ws3.Range("F12").Value = 0.7967082433047037
$ws3.Range("D15").Value = 31430.95
$ws3.Range("E15").Value = 23993.79316613378
$ws3.Range("F15").Value = 0.5670923887872029

# Column width tweaks on CUMPLIMIENTO MENSUAL (side effect of autofit).
# Excel's COM ColumnWidth units are offset from the stored OOXML <col width>
# by the fixed padding of 5/6 character (~0.8333) for this workbook's default
# font (Calibri 11) -- subtract it so the saved width lands on the exact
# integer value the diff expects.
$ws3.Columns.Item(4).ColumnWidth = 13 - 0.8333333
$ws3.Columns.Item(5).ColumnWidth = 23 - 0.8333333
